# Adds the "Tunniste" column header to G1 and updates the print/page setup,
# matching the commit "Added _TULOS to result file name AND added Tunniste-field"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Tunniste" field/header in column G (previously empty, style already set)
$ws.Range("G1").Value = "Tunniste"

# Excel moved the active selection to G1 after the edit
$ws.Range("G1").Select()

# Page setup was (re)written: A4 paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
